# Injuries_Master_Clubs.xlsx update (2025-11-15 refresh run)
# Sheet "snapshot": two players' injuries (Самсонов Илья, Хёфенмайер Ноэль) replace the
# rows formerly occupied by players who returned (Мачулин Василий, Хомченко Павел);
# the row for Кузьмин Глеб (also returned) is removed outright, shifting the remaining
# Салават Юлаев / Торпедо / ЦСКА / Драконы rows up by one. All "scraped_at" values are
# refreshed to the new run's timestamps.
# Sheet "returned": logs the 3 players who came off the injury list.
# Sheet "new_injured": logs the 2 players newly added to the injury list.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) snapshot sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("snapshot")

# Row 28: Мачулин Василий (СОЧ, защитник) -> Самсонов Илья (СОЧ, вратарь)
$ws.Range("D28").Value = "Самсонов Илья"
$ws.Range("E28").Value = "35"
$ws.Range("F28").Value = "вратарь"
$ws.Range("G28").Value = "21010"
$ws.Range("H28").Value = "1369_СОЧ_самсоновилья"

# Row 29: Хомченко Павел (СОЧ, вратарь) -> Хёфенмайер Ноэль (СОЧ, защитник)
$ws.Range("D29").Value = "Хёфенмайер Ноэль"
$ws.Range("E29").Value = "22"
$ws.Range("F29").Value = "защитник"
$ws.Range("G29").Value = "44847"
$ws.Range("H29").Value = "1369_СОЧ_хефенмайерноэль"

# Row 35: Кузьмин Глеб (СЮЛ) returned to play and is dropped from the list entirely;
# deleting the row shifts every following row (Салават Юлаев.. Драконы) up by one.
$ws.Rows.Item(35).Delete()

# Refresh "scraped_at" (column K) for every remaining data row with this run's timestamps.
$kTimestamps = @{
    2 = "2025-11-15T07:03:08.050729+00:00"
    3 = "2025-11-15T07:03:08.050772+00:00"
    4 = "2025-11-15T07:03:08.050795+00:00"
    5 = "2025-11-15T07:03:10.672959+00:00"
    6 = "2025-11-15T07:03:10.672990+00:00"
    7 = "2025-11-15T07:03:13.430032+00:00"
    8 = "2025-11-15T07:03:16.081522+00:00"
    9 = "2025-11-15T07:03:18.376428+00:00"
    10 = "2025-11-15T07:03:18.376459+00:00"
    11 = "2025-11-15T07:03:23.456790+00:00"
    12 = "2025-11-15T07:03:25.720710+00:00"
    13 = "2025-11-15T07:03:28.025406+00:00"
    14 = "2025-11-15T07:03:28.025438+00:00"
    15 = "2025-11-15T07:03:28.025459+00:00"
    16 = "2025-11-15T07:03:30.252069+00:00"
    17 = "2025-11-15T07:03:32.480356+00:00"
    18 = "2025-11-15T07:03:32.480392+00:00"
    19 = "2025-11-15T07:03:35.247095+00:00"
    20 = "2025-11-15T07:03:37.518757+00:00"
    21 = "2025-11-15T07:03:37.518792+00:00"
    22 = "2025-11-15T07:03:37.518816+00:00"
    23 = "2025-11-15T07:03:37.518840+00:00"
    24 = "2025-11-15T07:03:37.518862+00:00"
    25 = "2025-11-15T07:03:40.286678+00:00"
    26 = "2025-11-15T07:03:40.286708+00:00"
    27 = "2025-11-15T07:03:42.985477+00:00"
    28 = "2025-11-15T07:03:42.985507+00:00"
    29 = "2025-11-15T07:03:42.985527+00:00"
    30 = "2025-11-15T07:03:45.672530+00:00"
    31 = "2025-11-15T07:03:45.672557+00:00"
    32 = "2025-11-15T07:03:47.966011+00:00"
    33 = "2025-11-15T07:03:47.966048+00:00"
    34 = "2025-11-15T07:03:47.966070+00:00"
    35 = "2025-11-15T07:03:47.966092+00:00"
    36 = "2025-11-15T07:03:47.966112+00:00"
    37 = "2025-11-15T07:03:50.359413+00:00"
    38 = "2025-11-15T07:03:50.359445+00:00"
    39 = "2025-11-15T07:03:55.182226+00:00"
    40 = "2025-11-15T07:03:55.182258+00:00"
    41 = "2025-11-15T07:03:55.182277+00:00"
    42 = "2025-11-15T07:03:55.182295+00:00"
    43 = "2025-11-15T07:03:58.257072+00:00"
    44 = "2025-11-15T07:03:58.257104+00:00"
}

foreach ($rowNum in $kTimestamps.Keys) {
    $ws.Cells.Item($rowNum, 11).Value = $kTimestamps[$rowNum]
}

# ---------------------------------------------------------------------------
# 2) returned sheet - players who came off the injury list this run
# ---------------------------------------------------------------------------
$wsReturned = $wb.Worksheets.Item("returned")

# Format the "changed_day" column as text first so the ISO date string
# ("2025-11-15") is stored literally instead of being auto-converted to a
# date serial number.
$wsReturned.Range("G2:G4").NumberFormat = "@"

$wsReturned.Range("A2").Value = "СОЧ"
$wsReturned.Range("B2").Value = "ХК Сочи"
$wsReturned.Range("C2").Value = "Мачулин Василий"
$wsReturned.Range("D2").Value = "1369_СОЧ_мачулинвасилий"
$wsReturned.Range("E2").Value = "RETURN"
$wsReturned.Range("F2").Value = "2025-11-15T15:03:58.764737+08:00"
$wsReturned.Range("G2").Value = "2025-11-15"

$wsReturned.Range("A3").Value = "СОЧ"
$wsReturned.Range("B3").Value = "ХК Сочи"
$wsReturned.Range("C3").Value = "Хомченко Павел"
$wsReturned.Range("D3").Value = "1369_СОЧ_хомченкопавел"
$wsReturned.Range("E3").Value = "RETURN"
$wsReturned.Range("F3").Value = "2025-11-15T15:03:58.764737+08:00"
$wsReturned.Range("G3").Value = "2025-11-15"

$wsReturned.Range("A4").Value = "СЮЛ"
$wsReturned.Range("B4").Value = "Салават Юлаев"
$wsReturned.Range("C4").Value = "Кузьмин Глеб"
$wsReturned.Range("D4").Value = "1369_СЮЛ_кузьминглеб"
$wsReturned.Range("E4").Value = "RETURN"
$wsReturned.Range("F4").Value = "2025-11-15T15:03:58.764737+08:00"
$wsReturned.Range("G4").Value = "2025-11-15"

# ---------------------------------------------------------------------------
# 3) new_injured sheet - players newly added to the injury list this run
# ---------------------------------------------------------------------------
$wsNewInjured = $wb.Worksheets.Item("new_injured")

# Same text-formatting guard as above for the "changed_day" column.
$wsNewInjured.Range("G2:G3").NumberFormat = "@"

$wsNewInjured.Range("A2").Value = "СОЧ"
$wsNewInjured.Range("B2").Value = "ХК Сочи"
$wsNewInjured.Range("C2").Value = "Самсонов Илья"
$wsNewInjured.Range("D2").Value = "1369_СОЧ_самсоновилья"
$wsNewInjured.Range("E2").Value = "INJURED_NEW"
$wsNewInjured.Range("F2").Value = "2025-11-15T15:03:58.764737+08:00"
$wsNewInjured.Range("G2").Value = "2025-11-15"

$wsNewInjured.Range("A3").Value = "СОЧ"
$wsNewInjured.Range("B3").Value = "ХК Сочи"
$wsNewInjured.Range("C3").Value = "Хёфенмайер Ноэль"
$wsNewInjured.Range("D3").Value = "1369_СОЧ_хефенмайерноэль"
$wsNewInjured.Range("E3").Value = "INJURED_NEW"
$wsNewInjured.Range("F3").Value = "2025-11-15T15:03:58.764737+08:00"
$wsNewInjured.Range("G3").Value = "2025-11-15"
